$wb = $excel.ActiveWorkbook

# ALC!row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 121.76923
$ws.Range("I9").Value = 138.44444
$ws.Range("J9").Value = 84.25
$ws.Range("K9").Value = 138.44444
$ws.Range("L9").Value = 84.25
$ws.Range("M9").Value = 30.55556000000001
$ws.Range("N9").Value = -422.25

# ALC!row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 639.1923
$ws.Range("I18").Value = 652.76
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 652.76
$ws.Range("L18").Value = 300
$ws.Range("M18").Value = -368.76
$ws.Range("N18").Value = -868

# ALC!row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1025.6666
$ws.Range("I19").Value = 1223
$ws.Range("J19").Value = 965.6087
$ws.Range("K19").Value = 1223
$ws.Range("L19").Value = 965.6087
$ws.Range("M19").Value = -1048
$ws.Range("N19").Value = -1315.6087

# ALC!row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1308.9333
$ws.Range("I33").Value = 123.77778
$ws.Range("J33").Value = 3086.6667
$ws.Range("K33").Value = 123.77778
$ws.Range("L33").Value = 3086.6667
$ws.Range("M33").Value = 105.22222
$ws.Range("N33").Value = -3544.6667

# ALC!row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1063.1538
$ws.Range("I98").Value = 1015.7059
$ws.Range("J98").Value = 1152.7778
$ws.Range("K98").Value = 1015.7059
$ws.Range("L98").Value = 1152.7778
$ws.Range("M98").Value = 482.2941
$ws.Range("N98").Value = -4148.7778

# ALC!row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 571.2857
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 500
$ws.Range("L107").Value = 999
$ws.Range("M107").Value = 1420
$ws.Range("N107").Value = -4839

# ALC!row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 879.0714
$ws.Range("I112").Value = 600
$ws.Range("J112").Value = 990.7
$ws.Range("K112").Value = 1800
$ws.Range("L112").Value = 2972.1
$ws.Range("M112").Value = -692
$ws.Range("N112").Value = -5188.1

# ALC!row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1063.1538
$ws.Range("I122").Value = 1015.7059
$ws.Range("J122").Value = 1152.7778
$ws.Range("K122").Value = 3047.1177
$ws.Range("L122").Value = 3458.3334
$ws.Range("M122").Value = -597.1177000000002
$ws.Range("N122").Value = -8358.3334

# ALC!row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 557.51514
$ws.Range("I135").Value = 533.0968
$ws.Range("J135").Value = 936
$ws.Range("K135").Value = 4797.8712
$ws.Range("L135").Value = 8424
$ws.Range("M135").Value = -2262.8712
$ws.Range("N135").Value = -13494

# ALC!row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1399.75
$ws.Range("I137").Value = 923.7646999999999
$ws.Range("J137").Value = 1630.9429
$ws.Range("K137").Value = 2771.2941
$ws.Range("L137").Value = 4892.8287
$ws.Range("M137").Value = -221.2941000000001
$ws.Range("N137").Value = -9992.8287

# ARM!row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 27840542
$ws.Range("I110").Value = 35793704
$ws.Range("J110").Value = 4477.75
$ws.Range("K110").Value = 35793704
$ws.Range("L110").Value = 4477.75
$ws.Range("M110").Value = -35791659
$ws.Range("N110").Value = -8567.75

# BSM!row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1411.95
$ws.Range("I80").Value = 546.5
$ws.Range("J80").Value = 1988.9166
$ws.Range("K80").Value = 546.5
$ws.Range("L80").Value = 1988.9166
$ws.Range("M80").Value = 451.5
$ws.Range("N80").Value = -3984.9166

# BSM!row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 1411.95
$ws.Range("I83").Value = 546.5
$ws.Range("J83").Value = 1988.9166
$ws.Range("K83").Value = 2732.5
$ws.Range("L83").Value = 9944.583000000001
$ws.Range("M83").Value = 2259.5
$ws.Range("N83").Value = -19928.583

# CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17957.87
$ws.Range("I31").Value = 35520.516
$ws.Range("J31").Value = 2041.7188
$ws.Range("K31").Value = 35520.516
$ws.Range("L31").Value = 2041.7188
$ws.Range("M31").Value = -35225.516
$ws.Range("N31").Value = -2631.7188

# CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 17957.87
$ws.Range("I34").Value = 35520.516
$ws.Range("J34").Value = 2041.7188
$ws.Range("K34").Value = 35520.516
$ws.Range("L34").Value = 2041.7188
$ws.Range("M34").Value = -35318.516
$ws.Range("N34").Value = -2445.7188

# CRP!row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7528.3335
$ws.Range("I62").Value = 4005
$ws.Range("J62").Value = 9290
$ws.Range("K62").Value = 4005
$ws.Range("L62").Value = 9290
$ws.Range("M62").Value = -3381
$ws.Range("N62").Value = -10538

# CRP!row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 7528.3335
$ws.Range("I65").Value = 4005
$ws.Range("J65").Value = 9290
$ws.Range("K65").Value = 20025
$ws.Range("L65").Value = 46450
$ws.Range("M65").Value = -16905
$ws.Range("N65").Value = -52690

# CRP!row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13409.667
$ws.Range("I99").Value = 2319.75
$ws.Range("J99").Value = 22281.6
$ws.Range("K99").Value = 2319.75
$ws.Range("L99").Value = 22281.6
$ws.Range("M99").Value = -821.75
$ws.Range("N99").Value = -25277.6

# CRP!row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 998.9048
$ws.Range("I105").Value = 961.1875
$ws.Range("J105").Value = 1119.6
$ws.Range("K105").Value = 961.1875
$ws.Range("L105").Value = 1119.6
$ws.Range("M105").Value = 785.8125
$ws.Range("N105").Value = -4613.6

# CRP!row 112
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 43000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 43000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 43000
$ws.Range("N112").Value = -45954

# CRP!row 120
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 29800
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 29800
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 29800
$ws.Range("N120").Value = -37058

# CRP!row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2532.1304
$ws.Range("I122").Value = 2886.077
$ws.Range("J122").Value = 2072
$ws.Range("K122").Value = 8658.231
$ws.Range("L122").Value = 6216
$ws.Range("M122").Value = -6208.231
$ws.Range("N122").Value = -11116

# CRP!row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 13409.667
$ws.Range("I126").Value = 2319.75
$ws.Range("J126").Value = 22281.6
$ws.Range("K126").Value = 6959.25
$ws.Range("L126").Value = 66844.79999999999
$ws.Range("M126").Value = -4489.25
$ws.Range("N126").Value = -71784.79999999999

# CUL!row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 41.72222
$ws.Range("I12").Value = 10.75
$ws.Range("J12").Value = 50.57143
$ws.Range("K12").Value = 32.25
$ws.Range("L12").Value = 151.71429
$ws.Range("M12").Value = 140.75
$ws.Range("N12").Value = -497.71429

# CUL!row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 33.666668
$ws.Range("I14").Value = 33.666668
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 101.000004
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 71.999996

# CUL!row 103
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 117.666664
$ws.Range("I103").Value = 117.666664
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 352.999992
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = 526.000008
$ws.Range("N103").ClearContents()

# CUL!row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1430.08
$ws.Range("I131").Value = 318.75
$ws.Range("J131").Value = 1526.7174
$ws.Range("K131").Value = 956.25
$ws.Range("L131").Value = 4580.1522
$ws.Range("M131").Value = 4083.75
$ws.Range("N131").Value = -14660.1522

# CUL!row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2615
$ws.Range("I132").Value = 2600
$ws.Range("J132").Value = 2617.8572
$ws.Range("K132").Value = 23400
$ws.Range("L132").Value = 23560.7148
$ws.Range("M132").Value = -20870
$ws.Range("N132").Value = -28620.7148

# GSM!row 123
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 9726
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 9726
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 9726
$ws.Range("N123").Value = -14626

# GSM!row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3207
$ws.Range("I126").Value = 4336.2
$ws.Range("J126").Value = 2077.8
$ws.Range("K126").Value = 13008.6
$ws.Range("L126").Value = 6233.400000000001
$ws.Range("M126").Value = -10538.6
$ws.Range("N126").Value = -11173.4

# LTW!row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 103640.4
$ws.Range("I40").Value = 128238.125
$ws.Range("J40").Value = 5249.5
$ws.Range("K40").Value = 128238.125
$ws.Range("L40").Value = 5249.5
$ws.Range("M40").Value = -128102.125
$ws.Range("N40").Value = -5521.5

# LTW!row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3625
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3625
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3625
$ws.Range("N46").Value = -4001
$ws.Range("M46").ClearContents()

# LTW!row 131
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H131").Value = 37563
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 37563
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 37563
$ws.Range("N131").Value = -47643

# WVR!row 14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 11080
$ws.Range("I14").Value = 866.6667
$ws.Range("J14").Value = 26400
$ws.Range("K14").Value = 866.6667
$ws.Range("L14").Value = 26400
$ws.Range("M14").Value = -698.6667
$ws.Range("N14").Value = -26736

# WVR!row 20
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 3380
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 3380
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 3380
$ws.Range("N20").Value = -3860

# WVR!row 30
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 5745.3076
$ws.Range("I30").Value = 30994.5
$ws.Range("J30").Value = 1154.5454
$ws.Range("K30").Value = 30994.5
$ws.Range("L30").Value = 1154.5454
$ws.Range("M30").Value = -30887.5
$ws.Range("N30").Value = -1368.5454

# WVR!row 49
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 6808.857
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 6808.857
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 6808.857
$ws.Range("N49").Value = -7268.857

# WVR!row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1169
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1169
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 3507
$ws.Range("N122").Value = -8407
$ws.Range("M122").ClearContents()

# WVR!row 129
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 39990
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 39990
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 39990
$ws.Range("N129").Value = -49990
